# Excel COM-interop edit script
# Adds 8 rows of invoice/service-line data (two recurring 4-line "documents")
# to the worksheet, removes the two leftover styled-but-empty rows (10-11),
# and brings sheet view / column widths in line with the final layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Populate rows 2-9 with invoice line data ---
# Row 2
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "№ 514"
$ws.Range("C2").Value = "от 04 декабря 2021 r."
$ws.Range("D2").Value = "ООО ""КОМБИНАТ ИННОВАЦИОННЫХ ТЕХНОЛОГИЙ - MOHAPX"", ИНН"
$ws.Range("E2").Value = "ООО ""УК ""ГОРИЗОНТ"","
$ws.Range("F2").Value = "Работа мини погрузчика MUSTANG 2066 no договору комиссии NeK-01/21 от 01.04.2021, за период: Ноябрь 2021 г."
$ws.Range("G2").Value = 110
$ws.Range("H2").Value = "1 250,00"
$ws.Range("I2").Value = "137 500,00"

# Row 3
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "№ 514"
$ws.Range("C3").Value = "от 04 декабря 2021 r."
$ws.Range("D3").Value = "ООО ""КОМБИНАТ ИННОВАЦИОННЫХ ТЕХНОЛОГИЙ - MOHAPX"", ИНН"
$ws.Range("E3").Value = "ООО ""УК ""ГОРИЗОНТ"","
$ws.Range("F3").Value = "|Доставка минипогрузчика no договору комиссии №К-01/21 от 01.04.2021, за период: Ноябрь 2021 г."
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = "5 000,00"
$ws.Range("I3").Value = "5 000,00"

# Row 4
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "№ 3130"
$ws.Range("C4").Value = "30 ноября 2021 ."
$ws.Range("D4").Value = "ООО ""Комбинат Инновационных Технологий - МонАрх"""
$ws.Range("E4").Value = "ООО ""МонАрх и М"""
$ws.Range("F4").Value = "Размещение в гостинице Green Palace Vnukovo (Стандартный одноместный) Таланин А.В. 01.11.2021-13.11.2021 года"
$ws.Range("G4").Value = 12
$ws.Range("H4").Value = "3 000,00"
$ws.Range("I4").Value = "36 000,00"

# Row 5
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "№ 3130"
$ws.Range("C5").Value = "30 ноября 2021 ."
$ws.Range("D5").Value = "ООО ""Комбинат Инновационных Технологий - МонАрх"""
$ws.Range("E5").Value = "ООО ""МонАрх и М"""
$ws.Range("F5").Value = "Размещение в гостинице Green Palace Vnukovo (Стандартный одноместный) Ахычов Ибрагим 01.11.2021-04.12.2021 года"
$ws.Range("G5").Value = 33
$ws.Range("H5").Value = "2 100,00"
$ws.Range("I5").Value = "69 300,00"

# Row 6
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "№ 514"
$ws.Range("C6").Value = "от 04 декабря 2021 r."
$ws.Range("D6").Value = "ООО ""КОМБИНАТ ИННОВАЦИОННЫХ ТЕХНОЛОГИЙ - MOHAPX"", ИНН"
$ws.Range("E6").Value = "ООО ""УК ""ГОРИЗОНТ"","
$ws.Range("F6").Value = "Работа мини погрузчика MUSTANG 2066 no договору комиссии NeK-01/21 от 01.04.2021, за период: Ноябрь 2021 г."
$ws.Range("G6").Value = 110
$ws.Range("H6").Value = "1 250,00"
$ws.Range("I6").Value = "137 500,00"

# Row 7
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "№ 514"
$ws.Range("C7").Value = "от 04 декабря 2021 r."
$ws.Range("D7").Value = "ООО ""КОМБИНАТ ИННОВАЦИОННЫХ ТЕХНОЛОГИЙ - MOHAPX"", ИНН"
$ws.Range("E7").Value = "ООО ""УК ""ГОРИЗОНТ"","
$ws.Range("F7").Value = "|Доставка минипогрузчика no договору комиссии №К-01/21 от 01.04.2021, за период: Ноябрь 2021 г."
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = "5 000,00"
$ws.Range("I7").Value = "5 000,00"

# Row 8
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "№ 3130"
$ws.Range("C8").Value = "30 ноября 2021 ."
$ws.Range("D8").Value = "ООО ""Комбинат Инновационных Технологий - МонАрх"""
$ws.Range("E8").Value = "ООО ""МонАрх и М"""
$ws.Range("F8").Value = "Размещение в гостинице Green Palace Vnukovo (Стандартный одноместный) Таланин А.В. 01.11.2021-13.11.2021 года"
$ws.Range("G8").Value = 12
$ws.Range("H8").Value = "3 000,00"
$ws.Range("I8").Value = "36 000,00"

# Row 9
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "№ 3130"
$ws.Range("C9").Value = "30 ноября 2021 ."
$ws.Range("D9").Value = "ООО ""Комбинат Инновационных Технологий - МонАрх"""
$ws.Range("E9").Value = "ООО ""МонАрх и М"""
$ws.Range("F9").Value = "Размещение в гостинице Green Palace Vnukovo (Стандартный одноместный) Ахычов Ибрагим 01.11.2021-04.12.2021 года"
$ws.Range("G9").Value = 33
$ws.Range("H9").Value = "2 100,00"
$ws.Range("I9").Value = "69 300,00"

# --- Remove the two trailing placeholder rows that only carried a style ---
$ws.Rows("10:11").Delete()

# --- Column width / layout tidy-up to match the final sheet ---
$ws.Columns("D").ColumnWidth = 64
$ws.Columns("G").ColumnWidth = 6.25

# --- Selection / view state ---
[void]$ws.Range("A2:E5").Select()
